$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (from A1, which carries style index 1: bold/border/centered)
# onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every data row (2-59).
for ($r = 2; $r -le 59; $r++) {
    $ws.Cells.Item($r, 30).Value = 72
    $ws.Cells.Item($r, 31).Value = 90
    $ws.Cells.Item($r, 32).Value = 0
}
